# The commit swaps the two theme parts of the deck: the slide master's
# theme (ppt/theme/theme1.xml, currently the "Integral" / "Red Violet"
# color scheme used by every slide) takes on the stock "Office Theme"
# colors, while the notes master's theme (ppt/theme/theme2.xml) takes on
# the "Integral" / "Red Violet" colors that the slide master used to have.
#
# The only part of that swap reachable through the PowerPoint object
# model is the slide master's theme color scheme (every slide's Design
# ultimately points at it), so drive the recolor through
# Slide.ThemeColorScheme, which maps 1:1 onto the 12 DrawingML theme
# colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) in document
# order and persists straight into the clrScheme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Stock "Office" color scheme (RRGGBB -> OLE BGR long expected by .RGB)
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
